# Generate Report for Handoff
# This script updates the localization-status report:
#  - the handback status for 11620a51-59f4-49f4-9da4-f99daa341c5b.md moves from
#    "Handed back: in sync with en-US" to "Ready for handoff" (with new timestamps
#    and, where applicable, an explanatory error detail about a stale handback),
#  - and the two data rows on every sheet are re-ordered so that
#    34592a2d-0df0-442a-9a09-29b8a21fef50.md is listed first (row 2) and
#    11620a51-59f4-49f4-9da4-f99daa341c5b.md second (row 3).

$wb = $excel.ActiveWorkbook

$urlBase      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5576fcfe6ebb7366d9ed1aa8c7e2f600cae738d5/e2e"
$urlZhcn      = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/33fff5f3c1ac39f759fa3cc0285e1c00b180c2bc/e2e"
$urlDede      = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a462eac4ad6fd9c3dbc62b2e04c7f0c926d4c3ff/e2e"

$file11620 = "11620a51-59f4-49f4-9da4-f99daa341c5b.md"
$file34592 = "34592a2d-0df0-442a-9a09-29b8a21fef50.md"

$errorDetail = "The version of handback file is not the latest, current: $urlBase/$file11620, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edbb7857cf532753ef56ddff1a2188cdef3376a2/e2e/$file11620."

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $file34592
$wsOverview.Range("B2").Value = "e2e\$file34592"

$wsOverview.Range("A3").Value = $file11620
$wsOverview.Range("B3").Value = "e2e\$file11620"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 00:53:20"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$urlBase/$file11620", "", "", "e2e\$file34592") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$urlBase/$file34592", "", "", "e2e\$file11620") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhcn = $wb.Worksheets.Item("zh-cn")

$wsZhcn.Range("A2").Value = $file34592
$wsZhcn.Range("G2").Value = "34592a2d-0df0-442a-9a09-29b8a21fef50.33ba3e57dc62c9a01ed4d52789ca3bf53be9b056.zh-cn.xlf"
$wsZhcn.Range("I2").Value = $file34592
$wsZhcn.Range("J2").Value = "34592a2d-0df0-442a-9a09-29b8a21fef50.33ba3e57dc62c9a01ed4d52789ca3bf53be9b056.zh-cn.xlf"

$wsZhcn.Range("A3").Value = $file11620
$wsZhcn.Range("C3").Value = "Ready for handoff"
$wsZhcn.Range("G3").Value = "11620a51-59f4-49f4-9da4-f99daa341c5b.e8b555bcf62044172f0ba89631533206336c486d.zh-cn.xlf"
$wsZhcn.Range("H3").Value = "2016-09-04 00:53:15"
$wsZhcn.Range("I3").Value = $file11620
$wsZhcn.Range("J3").Value = "11620a51-59f4-49f4-9da4-f99daa341c5b.e8b555bcf62044172f0ba89631533206336c486d.zh-cn.xlf"
$wsZhcn.Range("P3").Value = $errorDetail

$wsZhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsZhcn.Hyperlinks.Delete()
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("A2"), "$urlBase/$file11620", "", "", $file34592) | Out-Null
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("I2"), "$urlZhcn/$file11620", "", "", $file34592) | Out-Null
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("A3"), "$urlBase/$file34592", "", "", $file11620) | Out-Null
$wsZhcn.Hyperlinks.Add($wsZhcn.Range("I3"), "$urlZhcn/$file34592", "", "", $file11620) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDede = $wb.Worksheets.Item("de-de")

$wsDede.Range("A2").Value = $file34592
$wsDede.Range("G2").Value = "34592a2d-0df0-442a-9a09-29b8a21fef50.33ba3e57dc62c9a01ed4d52789ca3bf53be9b056.de-de.xlf"
$wsDede.Range("I2").Value = $file34592
$wsDede.Range("J2").Value = "34592a2d-0df0-442a-9a09-29b8a21fef50.33ba3e57dc62c9a01ed4d52789ca3bf53be9b056.de-de.xlf"

$wsDede.Range("A3").Value = $file11620
$wsDede.Range("C3").Value = "Ready for handoff"
$wsDede.Range("G3").Value = "11620a51-59f4-49f4-9da4-f99daa341c5b.e8b555bcf62044172f0ba89631533206336c486d.de-de.xlf"
$wsDede.Range("H3").Value = "2016-09-04 00:53:20"
$wsDede.Range("I3").Value = $file11620
$wsDede.Range("J3").Value = "11620a51-59f4-49f4-9da4-f99daa341c5b.e8b555bcf62044172f0ba89631533206336c486d.de-de.xlf"
$wsDede.Range("P3").Value = $errorDetail

$wsDede.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDede.Hyperlinks.Delete()
$wsDede.Hyperlinks.Add($wsDede.Range("A2"), "$urlBase/$file11620", "", "", $file34592) | Out-Null
$wsDede.Hyperlinks.Add($wsDede.Range("I2"), "$urlDede/$file11620", "", "", $file34592) | Out-Null
$wsDede.Hyperlinks.Add($wsDede.Range("A3"), "$urlBase/$file34592", "", "", $file11620) | Out-Null
$wsDede.Hyperlinks.Add($wsDede.Range("I3"), "$urlDede/$file34592", "", "", $file11620) | Out-Null
